$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells remain text (not auto-converted to numbers)
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = "29.627.16"
$ws.Range("E2").Value = "  +0.41%  "
$ws.Range("D3").Value = "1.927.52"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.79%  "
$ws.Range("D5").Value = "326.75"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E6").Value = "  +0.82%  "
$ws.Range("D7").Value = "0.4828"
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("D8").Value = "0.4061"
$ws.Range("E8").Value = "  -0.48%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("D10").Value = "1.009"
$ws.Range("E10").Value = "  -0.96%  "
$ws.Range("D11").Value = "23.73"
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("D12").Value = "1.904.12"
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "6.070"
$ws.Range("E13").Value = "  +0.33%  "
$ws.Range("D14").Value = "7.289"
$ws.Range("E14").Value = "  +0.67%  "
$ws.Range("D15").Value = "91.46"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "0.06864"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").Value = "0.00001039"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "17.63"
$ws.Range("E19").Value = "  -0.90%  "
$ws.Range("D20").Value = "1.012"
$ws.Range("E20").Value = "  +0.67%  "
$ws.Range("D21").Value = "29.620.91"
$ws.Range("E21").Value = "  +0.32%  "
$ws.Range("D22").Value = "5.657"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("D23").Value = "11.95"
$ws.Range("E23").Value = "  +1.67%  "
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").Value = "2.154.50"
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "156.17"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "6.350"
$ws.Range("E27").Value = "  -4.38%  "
$ws.Range("E28").Value = "  -0.43%  "
$ws.Range("E29").Value = "  -2.15%  "
$ws.Range("D30").Value = "120.74"
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").Value = "1.003"
$ws.Range("E31").Value = "  -2.36%  "
$ws.Range("D32").Value = "0.09606"
$ws.Range("E32").Value = "  +0.53%  "
$ws.Range("D33").Value = "5.614"
$ws.Range("E33").Value = "  +1.47%  "
$ws.Range("D34").Value = "3.557"
$ws.Range("E34").Value = "  +0.07%  "
$ws.Range("D35").Value = "1.391"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("E36").Value = "  +6.70%  "
$ws.Range("D37").Value = "0.02279"
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "1.218"
$ws.Range("E38").Value = "  +2.56%  "
$ws.Range("E39").Value = "  -1.01%  "
$ws.Range("D40").Value = "10.74"
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("E41").Value = "  +0.61%  "
$ws.Range("D42").Value = "7.841"
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("D43").Value = "0.1844"
$ws.Range("E43").Value = "  -0.85%  "
$ws.Range("D44").Value = "2.487"
$ws.Range("E44").Value = "  +2.34%  "
$ws.Range("D45").Value = "1.281"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "0.07547"
$ws.Range("E47").Value = "  -1.13%  "
$ws.Range("D48").Value = "0.5551"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("D49").Value = "1.963"
$ws.Range("E49").Value = "  +0.11%  "
$ws.Range("D50").Value = "118.34"
$ws.Range("D51").Value = "2.432"
$ws.Range("E51").Value = "  +0.40%  "

# Reset style so no residual text-format style id lingers on the cells
$priceRange.Style = "Normal"
